$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 25; existing rows 25-39 shift down to 27-41.
$ws.Rows("25:26").Insert()

# New row 25 (weekly update entry)
$ws.Range("A25").Value = 9
$ws.Range("B25").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C25").Value = "Metropolitana"
$ws.Range("D25").Value = 44438
$ws.Range("E25").Value = 13
$ws.Range("F25").Value = 100114002
$ws.Range("G25").Value = "Camote"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 1240
$ws.Range("K25").Value = 12000
$ws.Range("L25").Value = 13000
$ws.Range("M25").Value = 12500
$ws.Range("N25").Value = "$/malla 18 kilos"
$ws.Range("O25").Value = "Perú"
$ws.Range("P25").Value = 694
$ws.Range("Q25").Value = 18
$ws.Range("R25").Value = "Hortaliza"

# New row 26 (weekly update entry)
$ws.Range("A26").Value = 9
$ws.Range("B26").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C26").Value = "Metropolitana"
$ws.Range("D26").Value = 44438
$ws.Range("E26").Value = 13
$ws.Range("F26").Value = 100114002
$ws.Range("G26").Value = "Camote"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Segunda"
$ws.Range("J26").Value = 700
$ws.Range("K26").Value = 10000
$ws.Range("L26").Value = 11000
$ws.Range("M26").Value = 10500
$ws.Range("N26").Value = "$/malla 18 kilos"
$ws.Range("O26").Value = "Perú"
$ws.Range("P26").Value = 583
$ws.Range("Q26").Value = 18
$ws.Range("R26").Value = "Hortaliza"
